# "chore: more monsters raaa"
#
# Fill in the "familiar" column (I) for the newly-added monster rows
# (52-74) with the value 20, formatted as text (numFmtId 49 / "@"),
# matching the style already used elsewhere in the sheet. Also update
# the sheet's scroll position / selection to reflect where the author
# was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

$range = $ws.Range("I52:I74")
$range.Value = 20
$range.NumberFormat = "@"

$excel.ActiveWindow.ScrollRow = 38
$excel.ActiveWindow.ScrollColumn = 1

$ws.Range("I52:I74").Select()
